$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-09"

# Update the header label cell (B1) that spells out the month/through-date text
$ws.Range("B1").Value = "October 2022 (through October 09)"

# --- New daily counts added across the "October" columns of every year ---

# Englewood (row 3)
$ws.Range("AP3").Value = 1

# Humboldt Park (row 4)
$ws.Range("L4").Value = 3

# Austin (row 5)
$ws.Range("L5").Value = 4
$ws.Range("V5").Value = 4

# North Lawndale (row 6)
$ws.Range("B6").Value = 1
$ws.Range("V6").Value = 4
$ws.Range("AZ6").Value = 1
$ws.Range("BJ6").Value = 1

# Grand Crossing (row 14)
$ws.Range("AZ14").Value = 2

# South Chicago (row 15)
$ws.Range("L15").Value = 2

# West Pullman (row 16)
$ws.Range("BT16").Value = 1

# Roseland (row 17)
$ws.Range("AZ17").Value = 1

# Auburn Gresham (row 23)
$ws.Range("B23").Value = 2

# Streeterville (row 33)
$ws.Range("AP33").Value = 1

# Loop (row 57)
$ws.Range("AZ57").Value = 2

# Lincoln Park (row 80)
$ws.Range("AP80").Value = 1

# North Center (row 86)
$ws.Range("AP86").Value = 2

# South Deering (row 94)
$ws.Range("V94").Value = 1

# West Elsdon (row 96)
$ws.Range("L96").Value = 1
